$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '56.738.22'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.332.02'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '515.58'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.68'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.336.50'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.37'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.16%  '
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.90'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.746.76'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '56.727.21'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.335.97'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '326.40'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('E21').Value = '  -0.74%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.61'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.84'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('E25').Value = '  +4.39%  '
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.40%  '
$ws.Range('E28').Value = '  +9.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '170.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.70'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.20'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.997'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.26'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.913'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  +2.79%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '38.33'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.88%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '142.88'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '277.29'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.20%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.16'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0505'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.01'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.98%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.54'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.94%  '
